# Format excel date cells as date for import/export
#
# The "transferred_at" column (A) currently stores its sample dates as the
# literal text "2020-01-01" (a shared string). Excel import/export expects a
# real date value there, so convert A2/A3 to numeric date serials and apply
# a yyyy-mm-dd date number format to them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 / A3: replace the text date with a true numeric date (serial 43831 =
# 2020-01-01) and format it as yyyy-mm-dd so it keeps displaying the same way.
$dateRange = $ws.Range("A2:A3")
$dateRange.Value = 43831
$dateRange.NumberFormat = "yyyy\-mm\-dd"

# Move the active selection (cosmetic, matches the saved workbook state).
$ws.Range("D25").Select()
